$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 291, shifting existing rows 291:305 down to 292:306
$ws.Rows.Item(291).Insert()

# Populate the new row 291 with the new price entry (same market/region/category
# context as the row that used to be there, but new date/price data).
$ws.Cells.Item(291, 1).Value = 11
$ws.Cells.Item(291, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(291, 3).Value = "Bíobío"
$ws.Cells.Item(291, 4).Value = 45041
$ws.Cells.Item(291, 5).Value = 8
$ws.Cells.Item(291, 6).Value = 100112040
$ws.Cells.Item(291, 7).Value = "Cilantro"
$ws.Cells.Item(291, 8).Value = "Sin especificar"
$ws.Cells.Item(291, 9).Value = "Primera"
$ws.Cells.Item(291, 10).Value = 180
$ws.Cells.Item(291, 11).Value = 7000
$ws.Cells.Item(291, 12).Value = 7500
$ws.Cells.Item(291, 13).Value = 7222
$ws.Cells.Item(291, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(291, 15).Value = "Región Metropolitana"
$ws.Cells.Item(291, 16).Value = 201
$ws.Cells.Item(291, 17).Value = 36
$ws.Cells.Item(291, 18).Value = "Hortaliza"
